# Update the "RES installed" sheet: increase installed RES capacity
# for nodes 6, 7 and 8 (rows 4-6, column C) from 15 to 30 MW.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RES installed")

$ws.Range("C4").Value = 30
$ws.Range("C5").Value = 30
$ws.Range("C6").Value = 30

# Reflect the selection left in the sheet at save time.
$ws.Range("D3").Select()

# Force a full recalculation so all dependent formulas
# (Main!B7, VLOOKUP-based Pg/Winter/Summer sheets, etc.) update.
$excel.CalculateFullRebuild()
